$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order-line rows to append below the existing header/data.
$data = @(
    @("INPSDR0120ITHACABAKERY", "Cup - Hot (12oz)",      "1",    "`$63.41",  "`$63.41"),
    @("H7658PC",                "Bag Trash - 38x58",     "1",    "`$38.11",  "`$38.11"),
    @("FIS508",                 "Bag - Wax (Sandwich)",  "0.17", "`$126.07", "`$21.43"),
    @("406020",                 "Spoon Soup - White MW", "1",    "`$10.43",  "`$10.43"),
    @("NPP406028",              "Fork - White HW",       "1",    "`$16.30",  "`$16.30")
)

$startRow = 3
$endRow = $startRow + $data.Length - 1

# Force the target range to text formatting first so that numeric-looking
# values ("1", "0.17", "406020") and currency-looking values ("$63.41")
# are written out as literal strings, not auto-converted numbers.
$targetRange = $ws.Range("A$startRow`:E$endRow")
$targetRange.NumberFormat = "@"

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Restore the default/Normal style so no extra number-format styling sticks
# to the new cells (matches the plain, unstyled cells used elsewhere).
$targetRange.Style = "Normal"
